# Update the people table with new/refreshed e-mail addresses and fill in
# the missing "interest" value for the last row, then re-create the
# mailto: hyperlinks so Excel re-applies its builtin "Hyperlink" cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing hyperlinks (and their special run-level formatting)
# before we touch the underlying cell values.
$ws.Hyperlinks.Delete()

# New e-mail addresses for the three people rows.
$ws.Range("C2").Value = "mvtukafke@emlhub.com"
$ws.Range("C3").Value = "bhlrhpdl@drope.ml"
$ws.Range("C4").Value = "ejwmipdl@drope.ml"

# Row 4 previously had no distinct interest value (it duplicated "nasa");
# give it its own value.
$ws.Range("D4").Value = "football"

# Re-insert the mailto hyperlinks on the (now updated) e-mail cells. Excel
# applies its builtin "Hyperlink" cell style automatically when doing this.
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:mvtukafke@emlhub.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:bhlrhpdl@drope.ml")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:ejwmipdl@drope.ml")
